$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row above row 234, shifting existing rows 234:239 down to 235:240
$ws.Rows.Item(234).Insert()

# Write the date as literal text (avoid Excel's automatic date-string parsing),
# then clear the temporary formatting so the cell matches its neighbors (no explicit style)
$ws.Cells.Item(234, 1).NumberFormat = "@"
$ws.Cells.Item(234, 1).Value = "12.01.2021"
$ws.Cells.Item(234, 1).ClearFormats()

# Populate the rest of the newly inserted row with the latest data point
$ws.Cells.Item(234, 2).Value = 38244
$ws.Cells.Item(234, 3).Value = 170099
$ws.Cells.Item(234, 4).Value = 1909
$ws.Cells.Item(234, 5).Value = 24461
$ws.Cells.Item(234, 6).Value = 11874
$ws.Cells.Item(234, 7).Value = 0
